$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.709.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.693.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4072"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.003"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08907"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.094"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001323"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.698.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07034"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.696.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.245"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.360"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.203"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.585"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08614"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.057"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.068"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2740"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.888"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09250"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.474"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7676"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.606"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7166"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.225"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.319"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "91.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07986"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
